# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6723
$wsExhibit.Range("F4").Value = 426
$wsExhibit.Range("F15").Value = 1456
$wsExhibit.Range("F17").Value = 3371
$wsExhibit.Range("F20").Value = 4
$wsExhibit.Range("F21").Value = 2010
$wsExhibit.Range("F22").Value = 121

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6723
$wsAll.Range("F4").Value = 426
$wsAll.Range("F16").Value = 1456
$wsAll.Range("F18").Value = 3371
$wsAll.Range("F21").Value = 4
$wsAll.Range("F22").Value = 2010
$wsAll.Range("F23").Value = 121
